$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.519.26'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.38%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.915.09'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.00%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.29'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.25%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.18%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4819'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +2.73%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2880'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.68%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06721'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.50%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '110.41'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.59%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.25'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +4.58%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.919.05'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.34%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07560'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.97%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.251'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.40%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6692'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.68%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '289.32'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -2.97%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.519.63'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.37%  '

$ws.Range("B18").Value = 'Dai'
$ws.Range("C18").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.001'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.19%  '

$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007589'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.44%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.91'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.04%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.167.37'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.17%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.482'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +4.85%  '

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.16%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.406'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +3.05%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.453'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.55%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.60'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.97%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.32'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -6.43%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.140'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.15%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1063'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.75%  '

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +3.24%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.162'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.03%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.023'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.97%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04987'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.44%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7295'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.03%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.133'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.67%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.20%  '

$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.739'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.05%  '

$ws.Range("B38").Value = 'Frax'
$ws.Range("C38").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.9994'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.23%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.15%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '110.70'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.90%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.010'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.34%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4412'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +3.79%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8640'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.87%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.897'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.43%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.22%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '68.41'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.33%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.322'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.79%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '48.72'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -5.23%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.308'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.10%  '

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +2.46%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.2508'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +3.38%  '
